# Update database for کیمیا-معدنی کیمیای زنجان گستران yearly income statement (rial)
# 1) Refresh company name label (B5)
# 2) Shift the 5-year rolling window of period headers (row 8) and published dates (row 9)
# 3) Shift all historical figures one column to the left and populate the new
#    rightmost column (H) with the newly published period's figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company name ---
$ws.Range("B5").Value = "کیمیا-ص. معدنی کیمیای زنجان گستران"

# --- Period headers (row 8) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Published dates (row 9) ---
$ws.Range("D9").Value = "1399-04-19 (13)"
$ws.Range("E9").Value = "1400-04-16 (12)"
$ws.Range("F9").Value = "1401-04-08 (9)"
$ws.Range("G9").Value = "1402-02-28 (8)"
# H9 looks like a bare date, so force text formatting to avoid Excel's automatic date conversion
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-28"

# --- Financial data (rows 11-27, columns D:H) ---
$ws.Range("D11").Value = 1185525
$ws.Range("E11").Value = 1505065
$ws.Range("F11").Value = 4246792
$ws.Range("G11").Value = 7078947
$ws.Range("H11").Value = 8846679
$ws.Range("D12").Value = -764099
$ws.Range("E12").Value = -939608
$ws.Range("F12").Value = -2759231
$ws.Range("G12").Value = -4991386
$ws.Range("H12").Value = -5785509
$ws.Range("D13").Value = 421426
$ws.Range("E13").Value = 565457
$ws.Range("F13").Value = 1487561
$ws.Range("G13").Value = 2087561
$ws.Range("H13").Value = 3061170
$ws.Range("D14").Value = -29566
$ws.Range("E14").Value = -33941
$ws.Range("F14").Value = -73858
$ws.Range("G14").Value = -317017
$ws.Range("H14").Value = -362437
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1672
$ws.Range("F16").Value = 10500
$ws.Range("G16").Value = 21000
$ws.Range("H16").Value = 923046
$ws.Range("D17").Value = 391860
$ws.Range("E17").Value = 533188
$ws.Range("F17").Value = 1424203
$ws.Range("G17").Value = 1791544
$ws.Range("H17").Value = 3621779
$ws.Range("D18").Value = -2745
$ws.Range("E18").Value = -16515
$ws.Range("F18").Value = -14844
$ws.Range("G18").Value = -3646
$ws.Range("H18").Value = -12575
$ws.Range("D19").Value = 23343
$ws.Range("E19").Value = 87426
$ws.Range("F19").Value = 148793
$ws.Range("G19").Value = 282249
$ws.Range("H19").Value = 186205
$ws.Range("D20").Value = 412458
$ws.Range("E20").Value = 604099
$ws.Range("F20").Value = 1558152
$ws.Range("G20").Value = 2070147
$ws.Range("H20").Value = 3795409
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("D22").Value = 412458
$ws.Range("E22").Value = 604099
$ws.Range("F22").Value = 1558152
$ws.Range("G22").Value = 2070147
$ws.Range("H22").Value = 3795409
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("D24").Value = 412458
$ws.Range("E24").Value = 604099
$ws.Range("F24").Value = 1558152
$ws.Range("G24").Value = 2070147
$ws.Range("H24").Value = 3795409
$ws.Range("D25").Value = 1031
$ws.Range("E25").Value = 575
$ws.Range("F25").Value = 1484
$ws.Range("G25").Value = 1972
$ws.Range("H25").Value = 1265
$ws.Range("D26").Value = 400000
$ws.Range("E26").Value = 1050000
$ws.Range("F26").Value = 1050000
$ws.Range("G26").Value = 1050000
$ws.Range("H26").Value = 3000000
$ws.Range("D27").Value = 137
$ws.Range("E27").Value = 201
$ws.Range("F27").Value = 519
$ws.Range("G27").Value = 690
$ws.Range("H27").Value = 1265
